$wb = $excel.ActiveWorkbook

$status = "Handed back: in sync with en-US"

# ---- Sheet "Overview" (unchanged data, just the status text) ----
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("B2").Value = $status
$ov.Range("C2").Value = $status
$ov.Range("B3").Value = $status
$ov.Range("C3").Value = $status

# ---- Helper data ----
$md1 = "2d4f67d6-7338-4056-b241-1ceadbdd4666.md"
$md2 = "ffffbe0ad6a2-7cf6-4cfa-b80a-ac88fe48ccb3.md"
$cfg = ".localization-config"

$md1Url   = "https://github.com/OpenLocalizationTest/oltest/blob/cb819d232096eeda08784cdbf0020b88e3d4d7da/e2e/2d4f67d6-7338-4056-b241-1ceadbdd4666.md"
$md2Url   = "https://github.com/OpenLocalizationTest/oltest/blob/cb819d232096eeda08784cdbf0020b88e3d4d7da/e2e/ffffbe0ad6a2-7cf6-4cfa-b80a-ac88fe48ccb3.md"
$cfgUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/cb819d232096eeda08784cdbf0020b88e3d4d7da/.localization-config"
$xlfZhUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1f43a31ab0322cd52fd80181e033f393aaf93fe2/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/2d4f67d6-7338-4056-b241-1ceadbdd4666.c6aeeeb8f5380877f1f78d926361dadc7c5dcd3c.zh-cn.xlf"
$xlfDeUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d99f35f0af0746b6d52e8a22a53d57feea1d99f3/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/2d4f67d6-7338-4056-b241-1ceadbdd4666.c6aeeeb8f5380877f1f78d926361dadc7c5dcd3c.de-de.xlf"

function Update-LangSheet($sheetName, $xlfName, $xlfUrl, $row2Date, $row4Date, $newHandbackDate) {
    $ws = $wb.Worksheets.Item($sheetName)

    # --- Row 2 (2d4f67d6...md) ---
    $ws.Range("A2").Value = $md1
    $ws.Range("B2").Value = $status
    $ws.Range("C2").Value = $xlfName
    $ws.Range("D2").Value = $row2Date

    $ws.Range("E2").Value = $md1
    $ws.Hyperlinks.Add($ws.Range("E2"), $md1Url, "", "", $md1) | Out-Null

    $ws.Range("F2").Value = $xlfName
    $ws.Hyperlinks.Add($ws.Range("F2"), $xlfUrl, "", "", $xlfName) | Out-Null

    $ws.Range("G2").Value = $newHandbackDate
    $ws.Range("H2").Value = "Include"

    # --- Row 3 (ffffbe0a...md) ---
    $ws.Range("A3").Value = $md2
    $ws.Range("B3").Value = $status
    $ws.Range("C3").Value = $xlfName
    $ws.Range("D3").Value = $row2Date

    $ws.Range("E3").Value = $md1
    $ws.Hyperlinks.Add($ws.Range("E3"), $md1Url, "", "", $md1) | Out-Null

    $ws.Range("F3").Value = $xlfName
    $ws.Hyperlinks.Add($ws.Range("F3"), $xlfUrl, "", "", $xlfName) | Out-Null

    $ws.Range("G3").Value = $newHandbackDate
    $ws.Range("H3").Value = "Include"

    # --- Row 4 (.localization-config) : untouched content, values stay the same ---
    $ws.Range("D4").Value = $row4Date
    $ws.Range("G4").Value = $row4Date
    $ws.Range("H4").Value = "Ignored"
}

Update-LangSheet "zh-cn" "2d4f67d6-7338-4056-b241-1ceadbdd4666.c6aeeeb8f5380877f1f78d926361dadc7c5dcd3c.zh-cn.xlf" $xlfZhUrl "2016-01-26 06:20:55" "0001-01-01 00:00:00" "2016-01-26 06:23:34"
Update-LangSheet "de-de" "2d4f67d6-7338-4056-b241-1ceadbdd4666.c6aeeeb8f5380877f1f78d926361dadc7c5dcd3c.de-de.xlf" $xlfDeUrl "2016-01-26 06:21:38" "0001-01-01 00:00:00" "2016-01-26 06:23:55"

Write-Host "done"
